$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (employeData -> list)
$ws.Name = "list"

# --- Row 1 (headers) ---
$ws.Range("A1").Value = 'firstName'
$ws.Range("B1").Value = 'middleName'
$ws.Range("C1").Value = 'lastName'
$ws.Range("D1").Value = 'userName'
$ws.Range("E1").Value = 'password'
$ws.Range("F1").ClearContents()

# --- Row 2 ---
$ws.Range("A2").Value = 'babar14'
$ws.Range("B2").Value = 'babar26'
$ws.Range("C2").Value = 'babar38'
$ws.Range("D2").Value = 'abcxyz83'
$ws.Range("E2").Value = 'aaAA1234eo05'
$ws.Range("F2").ClearContents()

# --- Row 3 ---
$ws.Range("A3").Value = 'pakin'
$ws.Range("B3").Value = 'mser'
$ws.Range("C3").Value = 'lsflsf'
$ws.Range("D3").Value = 'sk4i'
$ws.Range("E3").Value = 'lf34d$4#$'

# --- Row 4 ---
$ws.Range("A4").Value = 'injg4'
$ws.Range("B4").Value = 'sdf4rg'
$ws.Range("C4").Value = 'wf344'
$ws.Range("D4").Value = 'ret334'
$ws.Range("E4").Value = '5t4f4r5Frfg'

# Apply the new (black, vertically-centered) formatting to A3:F4.
# Build the exact style on a single cell first, then fan it out with a
# format-only paste so only one new style entry is added to the workbook.
$a3 = $ws.Range("A3")
$a3.VerticalAlignment = -4108
$a3.Font.Color = 0
$fmtRange = $ws.Range("A3:F4")
$a3.Copy()
$fmtRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false

# F3/F4 keep the new style but carry no value (matches the source row
# which only has data through column E from here on).
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()

# Column E widens to fit the longer password values.
$ws.Columns.Item(5).ColumnWidth = 12.3

# Selection moves to E4.
$ws.Range("E4").Select() | Out-Null
